$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header cell (H1) onto the new
# header cells I1 and J1 so they pick up the same cell style (bold font,
# thin border, centered alignment) without creating a new style entry.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I and J column data for rows 2 through 66
$iValues = @(8,9,9,9,6,8,9,9,6,7,8,6,8,6,7,6,9,5,5,8,7,5,11,5,4,7,6,8,9,8,9,6,6,7,9,8,9,9,7,6,7,7,7,9,8,9,9,9,10,8,7,9,8,8,9,9,8,6,8,8,8,7,6,7,5)
$jValues = @(9,9,9,9,6,8,9,9,7,7,8,6,9,7,7,7,9,6,6,8,7,5,12,5,5,7,7,8,9,8,9,6,7,8,9,9,9,10,7,7,7,8,7,9,10,9,9,9,10,8,7,9,8,9,9,9,8,6,8,8,8,7,6,7,5)

for ($r = 2; $r -le 66; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
